# This workbook's data rows (species observation records) were reshuffled:
# each destination row ends up holding the full contents that some other
# (source) row used to hold. Two independent permutation cycles are applied:
#   - rows 2-12 (cycle across each other)
#   - rows 14, 16, 17 (row 15 stays put)
# Strategy: snapshot every involved row's full A:AY contents first (so the
# permutation can be applied safely even though it touches overlapping rows),
# then write each snapshot into its new destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's original content moves here)
$mapping = @{
  2  = 7
  3  = 8
  4  = 5
  5  = 10
  6  = 2
  7  = 6
  8  = 12
  9  = 3
  10 = 11
  11 = 4
  12 = 9
  14 = 16
  15 = 15
  16 = 17
  17 = 14
}

$lastCol = "AY"

# Snapshot every source row's full contents before any writes happen.
$snapshots = @{}
foreach ($destRow in $mapping.Keys) {
  $srcRow = $mapping[$destRow]
  if (-not $snapshots.ContainsKey($srcRow)) {
    $addr = "A" + $srcRow + ":" + $lastCol + $srcRow
    $snapshots[$srcRow] = $ws.Range($addr).Value()
  }
}

# Now write the snapshots into their destination rows.
foreach ($destRow in $mapping.Keys) {
  $srcRow = $mapping[$destRow]
  $addr = "A" + $destRow + ":" + $lastCol + $destRow
  $ws.Range($addr).Value = $snapshots[$srcRow]
}
